$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the example well IDs with new values (A01 -> D01, A11 -> D11)
$ws.Range("A2").Value = "D01"
$ws.Range("A3").Value = "D11"

# Move the active selection to A4 (as if the user pressed Enter after editing A3)
$ws.Range("A4").Select()
